$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.501.30'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.71%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.181.82'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -4.07%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.48'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.71%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.608'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -6.54%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.191.30'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.69%  '
$ws.Range("E10").Value = '  -3.95%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.80'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.386'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.75%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.739.57'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.95%  '
$ws.Range("E14").Value = '  -2.36%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '64.552.79'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.59%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.44'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.77%  '
$ws.Range("E17").Value = '  -2.49%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.190.26'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.83%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '419.80'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.33%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.95'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.36'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.14'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.50%  '
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.56'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.91%  '
$ws.Range("B25").Value = 'LEO'
$ws.Range("C25").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.68'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.489'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.98%  '
$ws.Range("E28").Value = '  -6.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.79'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("E31").Value = '  -4.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.79'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.55%  '
$ws.Range("E33").Value = '  -0.09%  '
$ws.Range("E34").Value = '  -1.63%  '
$ws.Range("E35").Value = '  -3.73%  '
$ws.Range("E36").Value = '  -3.39%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '157.78'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.32%  '
$ws.Range("E38").Value = '  -5.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.734.59'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.03%  '
$ws.Range("E40").Value = '  -5.24%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '24.48'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.33%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.21'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.24'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.716'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.33%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.70'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.52%  '
$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0624'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.25%  '
$ws.Range("E47").Value = '  -2.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '21.76'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '294.55'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.14%  '
$ws.Range("B50").Value = 'FirstDigitalUSD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.998'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.21%  '
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.00'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -13.11%  '
